$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.116840481758118
$ws.Range("B1").Value = 3.649934768676758
$ws.Range("C1").Value = 4.4692063331604
$ws.Range("D1").Value = 1.877953171730042
$ws.Range("E1").Value = 1.304700016975403
